$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '28.809.75'
$cell.ClearFormats()
$ws.Range("E2").Value = '  +7.38%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.813.64'
$cell.ClearFormats()
$ws.Range("E3").Value = '  +5.08%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.9991'
$cell.ClearFormats()
$ws.Range("E4").Value = '  +0.15%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '250.86'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +4.01%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.4983'
$cell.ClearFormats()
$ws.Range("E7").Value = '  +2.48%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.2776'
$cell.ClearFormats()
$ws.Range("E8").Value = '  +7.29%  '

$ws.Range("E9").Value = '  +2.79%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.814.96'
$cell.ClearFormats()
$ws.Range("E10").Value = '  +5.13%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '16.78'
$cell.ClearFormats()
$ws.Range("E11").Value = '  +5.09%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.07162'
$cell.ClearFormats()
$ws.Range("E12").Value = '  +3.63%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.6495'
$cell.ClearFormats()
$ws.Range("E13").Value = '  +7.05%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '4.709'
$cell.ClearFormats()
$ws.Range("E14").Value = '  +5.27%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '82.02'
$cell.ClearFormats()
$ws.Range("E15").Value = '  +6.39%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '28.773.72'
$cell.ClearFormats()
$ws.Range("E16").Value = '  +8.13%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.9994'
$cell.ClearFormats()
$ws.Range("E17").Value = '  +0.12%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.000007387'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +3.12%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.9992'
$cell.ClearFormats()
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("E20").Value = '  +7.25%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '2.047.41'
$cell.ClearFormats()
$ws.Range("E21").Value = '  +4.96%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.611'
$cell.ClearFormats()
$ws.Range("E22").Value = '  +4.29%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '8.888'
$cell.ClearFormats()
$ws.Range("E23").Value = '  +3.67%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '5.355'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +5.43%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '143.96'
$cell.ClearFormats()
$ws.Range("E25").Value = '  +4.55%  '

$ws.Range("E26").Value = '  +4.73%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '117.78'
$cell.ClearFormats()
$ws.Range("E27").Value = '  +11.01%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.894'
$cell.ClearFormats()
$ws.Range("E28").Value = '  +6.89%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.396'
$cell.ClearFormats()
$ws.Range("E29").Value = '  +1.26%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.185'
$cell.ClearFormats()
$ws.Range("E30").Value = '  +6.46%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.08359'
$cell.ClearFormats()
$ws.Range("E31").Value = '  +4.84%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.852'
$cell.ClearFormats()
$ws.Range("E32").Value = '  +4.71%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.04964'
$cell.ClearFormats()
$ws.Range("E33").Value = '  +9.95%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.090'
$cell.ClearFormats()
$ws.Range("E34").Value = '  +7.82%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.6801'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +8.92%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '2.685'
$cell.ClearFormats()
$ws.Range("E36").Value = '  +3.29%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.739'
$cell.ClearFormats()
$ws.Range("E37").Value = '  +12.33%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.9693'
$cell.ClearFormats()
$ws.Range("E38").Value = '  +4.15%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.194'
$cell.ClearFormats()
$ws.Range("E39").Value = '  +7.89%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.01588'
$cell.ClearFormats()
$ws.Range("E40").Value = '  +6.15%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '6.007'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +6.42%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.9997'
$cell.ClearFormats()
$ws.Range("E42").Value = '  +0.17%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '101.46'
$cell.ClearFormats()

$ws.Range("E44").Value = '  +7.35%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '7.236'
$cell.ClearFormats()
$ws.Range("E45").Value = '  +5.74%  '

$ws.Range("E46").Value = '  +5.71%  '

$ws.Range("E47").Value = '  +1.94%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '8.170'
$cell.ClearFormats()
$ws.Range("E48").Value = '  +3.85%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '31.71'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +5.22%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.3652'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +8.36%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.310'
$cell.ClearFormats()
$ws.Range("E51").Value = '  +6.28%  '
